$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $ws.Range($cell).Style = "Normal"
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).Style = "Normal"
}

Set-TextValue "D2" "43.779.18"
Set-TextValue "E2" "  +0.00%  "
Set-TextValue "D3" "2.290.82"
Set-TextValue "E3" "  -0.18%  "
Set-TextValue "E4" "  +0.49%  "
Set-TextValue "D5" "113.46"
Set-TextValue "E5" "  +15.91%  "
Set-TextValue "D6" "268.27"
Set-TextValue "E6" "  -0.83%  "
Set-TextValue "D7" "0.626"
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "E8" "  +0.22%  "
Set-TextValue "D9" "0.616"
Set-TextValue "E9" "  +1.59%  "
Set-TextValue "D10" "48.59"
Set-TextValue "E10" "  +7.18%  "
Set-TextValue "D11" "0.0940"
Set-TextValue "E11" "  +0.53%  "
Set-TextValue "D12" "9.03"
Set-TextValue "E12" "  +14.28%  "
Set-TextValue "D13" "0.108"
Set-TextValue "E13" "  +0.52%  "
Set-TextValue "E14" "  -0.27%  "
Set-TextValue "D15" "2.635.17"
Set-TextValue "E15" "  -0.12%  "
Set-TextValue "D16" "0.867"
Set-TextValue "E16" "  +1.28%  "
Set-TextValue "D17" "2.290.52"
Set-TextValue "E17" "  +0.07%  "
Set-TextValue "D18" "43.677.78"
Set-TextValue "E18" "  -0.21%  "
Set-TextValue "D19" "0.0000109"
Set-TextValue "E19" "  -1.68%  "
Set-TextValue "D20" "7.01"
Set-TextValue "E20" "  +13.06%  "
Set-TextValue "D21" "72.24"
Set-TextValue "E21" "  -0.12%  "
Set-TextValue "E22" "  -0.68%  "
Set-TextValue "D23" "9.85"
Set-TextValue "E23" "  +7.87%  "
Set-TextValue "D24" "232.88"
Set-TextValue "E24" "  -0.15%  "
Set-TextValue "D25" "2.92"
Set-TextValue "E25" "  +3.16%  "
Set-TextValue "D26" "11.73"
Set-TextValue "E26" "  +3.83%  "
Set-TextValue "E27" "  -0.04%  "
Set-TextValue "D28" "43.15"
Set-TextValue "E28" "  +13.10%  "
Set-TextValue "D29" "3.91"
Set-TextValue "E29" "  +0.43%  "
Set-TextValue "E30" "  -2.23%  "
Set-TextValue "E31" "  +0.81%  "
Set-TextValue "D32" "174.17"
Set-TextValue "E32" "  -1.35%  "
Set-TextValue "D33" "0.0931"
Set-TextValue "E33" "  +4.01%  "
Set-TextValue "D34" "21.60"
Set-TextValue "E34" "  -0.89%  "
Set-TextValue "D35" "5.67"
Set-TextValue "E35" "  +4.41%  "
Set-TextValue "E36" "  +0.69%  "
Set-TextValue "E37" "  +0.91%  "
Set-TextValue "E38" "  +2.62%  "
Set-TextValue "E39" "  -1.32%  "
Set-TextValue "D40" "3.79"
Set-TextValue "E40" "  +7.41%  "
Set-TextValue "D41" "14.66"
Set-TextValue "E41" "  +20.47%  "
Set-TextValue "D42" "74.25"
Set-TextValue "E42" "  +14.63%  "
Set-TextValue "D43" "2.39"
Set-TextValue "E43" "  +2.55%  "
Set-TextValue "E44" "  +1.74%  "
Set-TextValue "D45" "6.33"
Set-TextValue "E45" "  +21.30%  "
Set-TextValue "E46" "  +0.23%  "
Set-TextValue "D47" "1.40"
Set-TextValue "E47" "  +1.86%  "
Set-TextValue "E48" "  -0.54%  "
Set-TextValue "D49" "103.00"
Set-TextValue "E49" "  +4.49%  "
Set-TextValue "E50" "  +3.95%  "
Set-TextValue "D51" "0.0999"
Set-TextValue "E51" "  -2.39%  "
